$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns: "_old" -> "_FV2210" and "_new" -> "_FV2304" ---
$headerMap = @{
    "Segmentname_old"          = "Segmentname_FV2210"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2210"
    "Segment_old"              = "Segment_FV2210"
    "Datenelement_old"         = "Datenelement_FV2210"
    "Segment ID_old"           = "Segment ID_FV2210"
    "Code_old"                 = "Code_FV2210"
    "Qualifier_old"            = "Qualifier_FV2210"
    "Beschreibung_old"         = "Beschreibung_FV2210"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2210"
    "Bedingung_old"            = "Bedingung_FV2210"
    "Segmentname_new"          = "Segmentname_FV2304"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2304"
    "Segment_new"              = "Segment_FV2304"
    "Datenelement_new"         = "Datenelement_FV2304"
    "Segment ID_new"           = "Segment ID_FV2304"
    "Code_new"                 = "Code_FV2304"
    "Qualifier_new"            = "Qualifier_FV2304"
    "Beschreibung_new"         = "Beschreibung_FV2304"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2304"
    "Bedingung_new"            = "Bedingung_FV2304"
}

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $current = $cell.Text
    if ($headerMap.ContainsKey($current)) {
        $cell.Value = $headerMap[$current]
    }
}

# --- Add a table (ListObject) covering the used range, with autofilter ---
$rng = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- Freeze the header row (top row) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
